$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray trailing placeholder row (if present) before shifting rows,
# so the row-insert below doesn't push it out past the valid row range.
$ws.Rows(1048576).Delete()

# Insert a new row above row 119 ("SE REQUIERE NUEVAMENTE" block), pushing
# the existing rows 119-156 down to 120-157. The new row inherits the
# formatting of row 118 automatically (style 5 on A/B, style 3 elsewhere).
$ws.Rows(119).Insert()

# Populate the newly inserted row 119 with the new requirement entry.
$ws.Range("A119").Value = "N10 - REQUIREMENT SPECIFICATION"
$ws.Range("B119").Value = "PARA QUE CONFIRME"
$ws.Range("D119").Value = "LAS CANTIDADES"

# Widen columns D and E to fit the newly added long text, and keep the
# following columns (now starting at F) at their original width.
$ws.Columns(4).ColumnWidth = 38.833333333333336
$ws.Columns(5).ColumnWidth = 26.666666666666668

# Reflect the new selection / scroll position used after the edit.
$ws.Range("D119").Select()
